$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 4589.4087
$ws.Range("I74").Value = 5627.8945
$ws.Range("J74").Value = 3393.5757
$ws.Range("K74").Value = 5627.8945
$ws.Range("L74").Value = 3393.5757
$ws.Range("M74").Value = -4691.8945
$ws.Range("N74").Value = -5265.575699999999

# Row 77
$ws.Range("H77").Value = 4589.4087
$ws.Range("I77").Value = 5627.8945
$ws.Range("J77").Value = 3393.5757
$ws.Range("K77").Value = 28139.4725
$ws.Range("L77").Value = 16967.8785
$ws.Range("M77").Value = -23459.4725
$ws.Range("N77").Value = -26327.8785

# Row 123
$ws.Range("H123").Value = 49520
$ws.Range("J123").Value = 49520
$ws.Range("L123").Value = 49520
$ws.Range("N123").Value = -59320

# Row 129
$ws.Range("H129").Value = 570.1875
$ws.Range("I129").Value = 393.9091
$ws.Range("J129").Value = 958
$ws.Range("K129").Value = 1181.7273
$ws.Range("L129").Value = 2874
$ws.Range("M129").Value = 3818.2727
$ws.Range("N129").Value = -12874

# Row 135
$ws.Range("H135").Value = 287190.72
$ws.Range("I135").Value = 287190.72
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 2584716.48
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -2582181.48
$ws.Range("N135").ClearContents()


$ws = $wb.Worksheets.Item("ARM")
# Row 37
$ws.Range("H37").Value = 21099.5
$ws.Range("J37").Value = 22786.066
$ws.Range("L37").Value = 22786.066
$ws.Range("N37").Value = -23332.066

# Row 55
$ws.Range("H55").Value = 34853
$ws.Range("J55").Value = 34853
$ws.Range("L55").Value = 34853
$ws.Range("N55").Value = -35483


$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 608.0769
$ws.Range("I94").Value = 519.84375
$ws.Range("K94").Value = 519.84375
$ws.Range("M94").Value = -68.84375

# Row 134
$ws.Range("H134").Value = 1532.6
$ws.Range("I134").Value = 1265.75
$ws.Range("J134").Value = 2600
$ws.Range("K134").Value = 3797.25
$ws.Range("L134").Value = 7800
$ws.Range("M134").Value = -1262.25
$ws.Range("N134").Value = -12870


$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 263.35294
$ws.Range("I22").Value = 252.07692
$ws.Range("K22").Value = 252.07692
$ws.Range("M22").Value = 97.92308

# Row 50
$ws.Range("H50").Value = 9096.333000000001
$ws.Range("J50").Value = 9096.333000000001
$ws.Range("L50").Value = 9096.333000000001
$ws.Range("N50").Value = -10346.333

# Row 60
$ws.Range("H60").Value = 26062.54
$ws.Range("J60").Value = 26062.54
$ws.Range("L60").Value = 26062.54
$ws.Range("N60").Value = -27084.54


$ws = $wb.Worksheets.Item("CUL")
# Row 87
$ws.Range("H87").Value = 12308.333
$ws.Range("I87").Value = 10140
$ws.Range("J87").Value = 13857.143
$ws.Range("K87").Value = 30420
$ws.Range("L87").Value = 41571.429
$ws.Range("M87").Value = -29172
$ws.Range("N87").Value = -44067.429

# Row 90
$ws.Range("H90").Value = 12308.333
$ws.Range("I90").Value = 10140
$ws.Range("J90").Value = 13857.143
$ws.Range("K90").Value = 91260
$ws.Range("L90").Value = 124714.287
$ws.Range("M90").Value = -85020
$ws.Range("N90").Value = -137194.287

# Row 92
$ws.Range("H92").Value = 838.125
$ws.Range("I92").Value = 607.1429000000001
$ws.Range("J92").Value = 1017.7778
$ws.Range("K92").Value = 1821.4287
$ws.Range("L92").Value = 3053.3334
$ws.Range("M92").Value = -573.4287000000002
$ws.Range("N92").Value = -5549.3334

# Row 96
$ws.Range("H96").Value = 4575
$ws.Range("J96").Value = 4575
$ws.Range("L96").Value = 13725
$ws.Range("N96").Value = -17843

# Row 108
$ws.Range("H108").Value = 1300
$ws.Range("I108").Value = 1300
$ws.Range("K108").Value = 3900
$ws.Range("M108").Value = -1020

# Row 109
$ws.Range("H109").Value = 6173
$ws.Range("I109").Value = 3513.5
$ws.Range("J109").Value = 6932.857
$ws.Range("K109").Value = 10540.5
$ws.Range("L109").Value = 20798.571
$ws.Range("M109").Value = -9500.5
$ws.Range("N109").Value = -22878.571

# Row 111
$ws.Range("H111").Value = 1113
$ws.Range("I111").Value = 1113
$ws.Range("K111").Value = 3339
$ws.Range("M111").Value = -272

# Row 126
$ws.Range("H126").Value = 1680.5264
$ws.Range("I126").Value = 1030
$ws.Range("J126").Value = 1716.6666
$ws.Range("K126").Value = 3090
$ws.Range("L126").Value = 5149.9998
$ws.Range("M126").Value = 1850
$ws.Range("N126").Value = -15029.9998

# Row 139
$ws.Range("H139").Value = 22505.55
$ws.Range("I139").Value = 1540.7142
$ws.Range("J139").Value = 50458.668
$ws.Range("K139").Value = 4622.142599999999
$ws.Range("L139").Value = 151376.004
$ws.Range("M139").Value = 517.8574000000008
$ws.Range("N139").Value = -161656.004


$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 12259.5
$ws.Range("J57").Value = 14661
$ws.Range("L57").Value = 14661
$ws.Range("N57").Value = -16301

# Row 62
$ws.Range("H62").Value = 31715
$ws.Range("J62").Value = 31715
$ws.Range("L62").Value = 31715
$ws.Range("N62").Value = -33087

# Row 65
$ws.Range("H65").Value = 31715
$ws.Range("J65").Value = 31715
$ws.Range("L65").Value = 95145
$ws.Range("N65").Value = -102009

# Row 133
$ws.Range("H133").Value = 24696
$ws.Range("J133").Value = 24696
$ws.Range("L133").Value = 24696
$ws.Range("N133").Value = -34816

# Row 138
$ws.Range("H138").Value = 66685.57000000001
$ws.Range("J138").Value = 66685.57000000001
$ws.Range("L138").Value = 66685.57000000001
$ws.Range("N138").Value = -76965.57000000001


$ws = $wb.Worksheets.Item("LTW")
# Row 63
$ws.Range("H63").Value = 31338.334
$ws.Range("J63").Value = 31338.334
$ws.Range("L63").Value = 31338.334
$ws.Range("N63").Value = -32836.334

# Row 66
$ws.Range("H66").Value = 31338.334
$ws.Range("J66").Value = 31338.334
$ws.Range("L66").Value = 94015.00199999999
$ws.Range("N66").Value = -101503.002

# Row 93
$ws.Range("H93").Value = 1820.7142
$ws.Range("I93").Value = 1739.1875
$ws.Range("J93").Value = 2081.6
$ws.Range("K93").Value = 1739.1875
$ws.Range("L93").Value = 2081.6
$ws.Range("M93").Value = -491.1875
$ws.Range("N93").Value = -4577.6

# Row 132
$ws.Range("H132").Value = 5810.96
$ws.Range("I132").Value = 5856.9473
$ws.Range("J132").Value = 5665.3335
$ws.Range("K132").Value = 17570.8419
$ws.Range("L132").Value = 16996.0005
$ws.Range("M132").Value = -15040.8419
$ws.Range("N132").Value = -22056.0005

# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()


$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 42143
$ws.Range("J46").Value = 42143
$ws.Range("L46").Value = 42143
$ws.Range("N46").Value = -42605

# Row 134
$ws.Range("H134").Value = 42143
$ws.Range("J134").Value = 42143
$ws.Range("L134").Value = 126429
$ws.Range("N134").Value = -131499

# Row 136
$ws.Range("H136").Value = 1681.55
$ws.Range("I136").Value = 1695.9412
$ws.Range("J136").Value = 1600
$ws.Range("K136").Value = 5087.8236
$ws.Range("L136").Value = 4800
$ws.Range("M136").Value = -2537.8236
$ws.Range("N136").Value = -9900

